$d = $word.ActiveDocument

# 1) Phone number -> formatted FR phone text
$d.Content.Find.Execute("06663240636", $true, $false, $false, $false, $false, $true, 1, $false, "Tél (FR) : +33 66 63 24 06 36", 2)

# 2) Insert "Mail : " right before the email hyperlink text (after the line break
#    that follows the phone number), so the contact block gains a "Mail :" label.
$rMail = $d.Content
$rMail.Find.Execute("alexandre.poitevin")
$mailInsertPoint = $d.Range($rMail.Start, $rMail.Start)
$mailInsertPoint.InsertBefore("Mail : ")

# 3) Skills line: add Flask to the intermediate skill set
$d.Content.Find.Execute(": SQL, Kivy", $true, $false, $false, $false, $false, $true, 1, $false, ": SQL, Kivy, Flask", 2)

# 4) Languages line: add native French
$d.Content.Find.Execute(": Anglais (bon niveau), Arabe littéraire (bon niveau)", $true, $false, $false, $false, $false, $true, 1, $false, ": Français (natif), Anglais (bon niveau), Arabe littéraire (bon niveau)", 2)
